$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.755.73'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '1.633.51'
$ws.Range("E3").Value = '  -0.20%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '
$ws.Range("E6").Value = '  -0.89%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.259'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  -3.46%  '
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("E12").Value = '  -0.20%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").Value = '1.858.28'
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.632.65'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("E15").Value = '  -0.27%  '
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.78'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").Value = '25.762.23'
$ws.Range("E18").Value = '  -0.44%  '
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("E20").Value = '  +1.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '194.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.54%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.94%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.53'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = '1.129.45'
$ws.Range("E37").Value = '  -0.38%  '
$ws.Range("E38").Value = '  -1.82%  '
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  -0.89%  '
$ws.Range("E41").Value = '  +0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.02%  '
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").Value = '1.768.24'
$ws.Range("E45").Value = '  -0.40%  '
$ws.Range("D46").Value = '0.0₆0110'
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.417'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.35%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.56'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
